$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Select()

# Add new rule row "Assign Alfresco Folder" to the RuleTable (row 23)
$ws.Range("B23").Value = "Assign Alfresco Folder"
$ws.Range("C23").Value = "container?.folder?.cmisFolderId == null"
$ws.Range("D23").Value = "setEcmFolderPath, '/Sites/acm/documentLibrary/Timesheets/' + dateFormat('yyyyMMdd') + '_' + `$acmTimesheet.getId()"

# Widen column D to fit the new, longer action text
$ws.Columns.Item(4).ColumnWidth = 111.34

# Match the author's final scroll position / selection
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D37").Select()

Write-Output "done"
